# Commit: "added HUQ051 to questionnaire"
#
# Inserts a new worksheet "HUQ051" into the questionnaire workbook right
# after "HUQ071" (i.e. before "OHQ033"), populates it with the
# Answer/Coefficient crosswalk table, autosizes column B, and leaves the
# selection state matching the authored edit. Also nudges the selection on
# "KIQ005" to A1:B1, matching an incidental selection change in the source
# edit.

$wb = $excel.ActiveWorkbook

# --- Minor selection change observed on the existing KIQ005 sheet -----
$kiq = $wb.Worksheets.Item("KIQ005")
$kiq.Range("A1:B1").Select() | Out-Null

# --- Insert the new HUQ051 sheet right after HUQ071 --------------------
$after = $wb.Worksheets.Item("HUQ071")
$new = $wb.Worksheets.Add($null, $after)
$new.Name = "HUQ051"

# --- Header row ----------------------------------------------------------
$new.Range("A1").Value = "Answer"
$new.Range("B1").Value = "Coefficient"

# Note: cell-write order below matters -- it reproduces the exact order in
# which new shared strings were introduced by the original author.
$new.Range("B13").Value = "."

$new.Range("A2").Value = "None"
$new.Range("B2").Value = 0

$new.Range("A3").Value = 1
$new.Range("B3").Value = 1

$new.Range("A4").Value = "2 to 3"
$new.Range("B4").Value = 2

$new.Range("A5").Value = "4 to 5"
$new.Range("B5").Value = 3

$new.Range("A6").Value = "6 to 7"
$new.Range("B6").Value = 4

$new.Range("A7").Value = "8 to 9"
$new.Range("B7").Value = 5

$new.Range("A8").Value = "10 to 12"
$new.Range("B8").Value = 6

$new.Range("A9").Value = "13 to 15"
$new.Range("B9").Value = 7

$new.Range("A10").Value = "16 or more"
$new.Range("B10").Value = 8

$new.Range("A11").Value = "Refused"
$new.Range("B11").Value = 77

$new.Range("A12").Value = "Don’t Know"
$new.Range("B12").Value = 99

$new.Range("A13").Value = "Missing"
# B13 already set to "." above.

# --- Column sizing + final selection, matching authored state ----------
$new.Columns("B:B").AutoFit() | Out-Null
$new.Range("I33").Select() | Out-Null
